$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data to append: each entry corresponds to columns B,C,D,E,F
# (A / id is computed as 52 + row index, matching ids 53..78)
$data = @(
    @("Metal Mario - Tennis", "Metal Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Super Mario Cereal", "Mario Cereal", "Others", "Kellogs", "Card"),
    @("Baby Mario - Soccer", "Baby Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Metal Mario - Soccer", "Metal Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Mario - Soccer", "Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Mario", "Mario", "Super Mario Bros.", "Super Mario", "Figure"),
    @("8-Bit Mario Modern Color", "Mario", "8-bit Mario", "Super Mario", "Figure"),
    @("Dr. Mario", "Mario", "Super Smash Bros.", "Super Mario", "Figure"),
    @("Baby Mario - Horse Racing", "Baby Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Mario - Cat", "Mario", "Super Mario Bros.", "Super Mario", "Figure"),
    @("Baby Mario - Golf", "Baby Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Mario - Wedding", "Mario", "Super Mario Bros.", "Super Mario", "Figure"),
    @("Metal Mario - Golf", "Metal Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Mario", "Mario", "Super Smash Bros.", "Super Mario", "Figure"),
    @("Mario - Tennis", "Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Baby Mario - Tennis", "Baby Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Mario - Gold Edition", "Mario", "Super Mario Bros.", "Super Mario", "Figure"),
    @("Mario - Power Up Band", "Mario", "Super Nintendo World", "Super Mario", "Band"),
    @("Metal Mario - Horse Racing", "Metal Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Baby Mario - Baseball", "Baby Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Mario - Golf", "Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Mario - Silver Edition", "Mario", "Super Mario Bros.", "Super Mario", "Figure"),
    @("Mario - Horse Racing", "Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("8-Bit Mario Classic Color", "Mario", "8-bit Mario", "Super Mario", "Figure"),
    @("Metal Mario - Baseball", "Metal Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card"),
    @("Mario - Baseball", "Mario", "Mario Sports Superstars", "Mario Sports Superstars", "Card")
)

$startRow = 54
$startId = 53

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $entry = $data[$i]

    $ws.Cells.Item($row, 1).Value = $startId + $i
    $ws.Cells.Item($row, 2).Value = $entry[0]
    $ws.Cells.Item($row, 3).Value = $entry[1]
    $ws.Cells.Item($row, 4).Value = $entry[2]
    $ws.Cells.Item($row, 5).Value = $entry[3]
    $ws.Cells.Item($row, 6).Value = $entry[4]
}
